# fix(publipostage): Try to solve Excel emoji problem
#
# Replace the status markers used in column A ("statut") with new
# markers:
#   old -> new
#   📘  -> ⚠️
#   📕  -> -3   (kept as TEXT, not a number)
#   📗  -> ✅

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$oldToNew = @{
    "📘" = "⚠️"
    "📕" = "-3"
    "📗" = "✅"
}

$usedRange = $ws.UsedRange
$rowCount = $usedRange.Rows.Count

for ($r = 1; $r -le $rowCount; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    $val = $cell.Value2
    if ($oldToNew.ContainsKey($val)) {
        $newVal = $oldToNew[$val]
        # "-3" looks numeric; force text storage so it round-trips as a
        # string (matches the other status markers, e.g. "⚠️"/"✅")
        # instead of being coerced into a number.
        if ($newVal -eq "-3") {
            $cell.NumberFormat = "@"
        }
        $cell.Value = $newVal
    }
}
